$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2899.5
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# Row 32
$ws.Range("H32").Value = 950
$ws.Range("I32").Value = 950
$ws.Range("K32").Value = 950
$ws.Range("M32").Value = -624
# Row 51
$ws.Range("H51").Value = 13312.375
$ws.Range("I51").Value = 6750
$ws.Range("K51").Value = 6750
$ws.Range("M51").Value = -6266
# Row 53
$ws.Range("H53").Value = 243.2
$ws.Range("I53").Value = 69.5
$ws.Range("J53").Value = 359
$ws.Range("K53").Value = 69.5
$ws.Range("L53").Value = 359
$ws.Range("M53").Value = 567.5
$ws.Range("N53").Value = -1633
# Row 100
$ws.Range("H100").Value = 3573.25
$ws.Range("J100").Value = 3447
$ws.Range("L100").Value = 3447
$ws.Range("N100").Value = -4529
# Row 107
$ws.Range("H107").Value = 602.913
$ws.Range("I107").Value = 618.8
$ws.Range("J107").Value = 497
$ws.Range("K107").Value = 618.8
$ws.Range("L107").Value = 497
$ws.Range("M107").Value = 1301.2
$ws.Range("N107").Value = -4337
# Row 129
$ws.Range("H129").Value = 1420.2858
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 2415.6
$ws.Range("I88").Value = 2749.2
$ws.Range("J88").Value = 2082
$ws.Range("K88").Value = 2749.2
$ws.Range("L88").Value = 2082
$ws.Range("M88").Value = -2343.2
$ws.Range("N88").Value = -2894
# Row 91
$ws.Range("H91").Value = 2415.6
$ws.Range("I91").Value = 2749.2
$ws.Range("J91").Value = 2082
$ws.Range("K91").Value = 2749.2
$ws.Range("L91").Value = 2082
$ws.Range("M91").Value = -1345.2
$ws.Range("N91").Value = -4890
# Row 122
$ws.Range("H122").Value = 1798.3334
$ws.Range("I122").Value = 1798.3334
$ws.Range("K122").Value = 5395.0002
$ws.Range("M122").Value = -2945.0002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value = 23783.166
$ws.Range("J88").Value = 23783.166
$ws.Range("L88").Value = 23783.166
$ws.Range("N88").Value = -24595.166
# Row 91
$ws.Range("H91").Value = 23783.166
$ws.Range("J91").Value = 23783.166
$ws.Range("L91").Value = 23783.166
$ws.Range("N91").Value = -26591.166

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 870.5714
$ws.Range("I22").Value = 899
$ws.Range("J22").Value = 766.3333
$ws.Range("K22").Value = 899
$ws.Range("L22").Value = 766.3333
$ws.Range("M22").Value = -549
$ws.Range("N22").Value = -1466.3333
# Row 88
$ws.Range("H88").Value = 23590.334
$ws.Range("J88").Value = 23590.334
$ws.Range("L88").Value = 23590.334
$ws.Range("N88").Value = -24402.334
# Row 91
$ws.Range("H91").Value = 23590.334
$ws.Range("J91").Value = 23590.334
$ws.Range("L91").Value = 23590.334
$ws.Range("N91").Value = -26398.334
# Row 99
$ws.Range("H99").Value = 1263575
$ws.Range("I99").Value = 850600
$ws.Range("K99").Value = 850600
$ws.Range("M99").Value = -849102
# Row 126
$ws.Range("H126").Value = 1263575
$ws.Range("I126").Value = 850600
$ws.Range("K126").Value = 2551800
$ws.Range("M126").Value = -2549330
# Row 132
$ws.Range("H132").Value = 1954.091
$ws.Range("I132").Value = 1954.091
$ws.Range("K132").Value = 5862.272999999999
$ws.Range("M132").Value = -3332.272999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 20.125
$ws.Range("I7").Value = 19.833334
$ws.Range("K7").Value = 59.500002
$ws.Range("M7").Value = 52.499998
# Row 107
$ws.Range("H107").Value = 444.66666
$ws.Range("I107").Value = 407.57144
$ws.Range("J107").Value = 477.125
$ws.Range("K107").Value = 1222.71432
$ws.Range("L107").Value = 1431.375
$ws.Range("M107").Value = 697.28568
$ws.Range("N107").Value = -5271.375
# Row 132
$ws.Range("H132").Value = 589
$ws.Range("I132").Value = 230
$ws.Range("J132").Value = 1127.5
$ws.Range("K132").Value = 2070
$ws.Range("L132").Value = 10147.5
$ws.Range("M132").Value = 460
$ws.Range("N132").Value = -15207.5
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 14664.619
$ws.Range("I7").Value = 14553.444
$ws.Range("K7").Value = 14553.444
$ws.Range("M7").Value = -14441.444
# Row 16
$ws.Range("H16").Value = 1191.1666
$ws.Range("I16").Value = 1029.4
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1029.4
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -859.4000000000001
$ws.Range("N16").Value = -2340
# Row 68
$ws.Range("H68").Value = 3549.75
$ws.Range("I68").Value = 3386.2856
$ws.Range("J68").Value = 3778.6
$ws.Range("K68").Value = 3386.2856
$ws.Range("L68").Value = 3778.6
$ws.Range("M68").Value = -2637.2856
$ws.Range("N68").Value = -5276.6
# Row 71
$ws.Range("H71").Value = 3549.75
$ws.Range("I71").Value = 3386.2856
$ws.Range("J71").Value = 3778.6
$ws.Range("K71").Value = 16931.428
$ws.Range("L71").Value = 18893
$ws.Range("M71").Value = -13187.428
$ws.Range("N71").Value = -26381
# Row 82
$ws.Range("H82").Value = 4562.5
$ws.Range("I82").Value = 4166.6665
$ws.Range("J82").Value = 4800
$ws.Range("K82").Value = 4166.6665
$ws.Range("L82").Value = 4800
$ws.Range("M82").Value = -3805.6665
$ws.Range("N82").Value = -5522
# Row 85
$ws.Range("H85").Value = 4562.5
$ws.Range("I85").Value = 4166.6665
$ws.Range("J85").Value = 4800
$ws.Range("K85").Value = 4166.6665
$ws.Range("L85").Value = 4800
$ws.Range("M85").Value = -2918.6665
$ws.Range("N85").Value = -7296
# Row 126
$ws.Range("H126").Value = 14664.619
$ws.Range("I126").Value = 14553.444
$ws.Range("K126").Value = 43660.33199999999
$ws.Range("M126").Value = -41190.33199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4883.6
$ws.Range("I62").Value = 4854.5
$ws.Range("K62").Value = 4854.5
$ws.Range("M62").Value = -4230.5
# Row 65
$ws.Range("H65").Value = 4883.6
$ws.Range("I65").Value = 4854.5
$ws.Range("K65").Value = 24272.5
$ws.Range("M65").Value = -21152.5
# Row 81
$ws.Range("H81").Value = 2082.2
$ws.Range("I81").Value = 2082.2
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4164.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -3103.4
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 2082.2
$ws.Range("I84").Value = 2082.2
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 20822
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -15518
$ws.Range("N84").ClearContents()
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
